$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): update B1:D1 text, add E1:G1 ---
$ws.Range("A1").Value = "Core Attribute"
$ws.Range("B1").Value = "test1.csv-COL2"
$ws.Range("C1").Value = "test1.csv-COL3"
$ws.Range("D1").Value = "test1.csv-COL4"
$ws.Range("E1").Value = "http://dbpedia.org/ontology/percentage"
$ws.Range("F1").Value = "test1.csv-COL6"
$ws.Range("G1").Value = "http://dbpedia.org/ontology/populationTotal"

# Copy header style (bold, border, centered) from A1 onto the new E1:G1 header cells
$ws.Cells.Item(1, 1).Copy()
$ws.Range("E1:G1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 2 ---
$ws.Range("A2").Value = "http://dbpedia.org/resource/Charny-sur-Meuse"
$ws.Range("B2").Value = "http://dbpedia.org/resource/NAMUR"
$ws.Range("C2").Value = "'92075"
$ws.Range("D2").Value = "'5101"
$ws.Range("E2").Value = "'16"
$ws.Range("F2").Value = "http://dbpedia.org/resource/LM"
$ws.Range("G2").Value = "'476"

# --- Row 3 ---
$ws.Range("A3").Value = "http://dbpedia.org/resource/Flawinne"
$ws.Range("B3").Value = "http://dbpedia.org/resource/NAMUR"
$ws.Range("C3").Value = "'92043"
$ws.Range("D3").Value = "'5020"
$ws.Range("E3").Value = "'71"
$ws.Range("F3").Value = "http://dbpedia.org/resource/FW"
$ws.Range("G3").Value = "'4491"

# --- Row 4 ---
$ws.Range("A4").Value = "http://dbpedia.org/resource/Terp"
$ws.Range("B4").Value = "http://dbpedia.org/resource/NAMUR"
$ws.Range("C4").Value = "'92136"
$ws.Range("D4").Value = "'5100"
$ws.Range("E4").Value = "'54"
$ws.Range("F4").Value = "http://dbpedia.org/resource/WD"
$ws.Range("G4").Value = "'1845"

# --- Row 5 ---
$ws.Range("A5").Value = "Naninne"
$ws.Range("B5").Value = "http://dbpedia.org/resource/NAMUR"
$ws.Range("C5").Value = "'92095"
$ws.Range("D5").Value = "'5100"
$ws.Range("E5").Value = "'45"
$ws.Range("F5").Value = "http://dbpedia.org/resource/NN"
$ws.Range("G5").Value = "'1606"

# --- Row 6 ---
$ws.Range("A6").Value = "http://dbpedia.org/resource/Dave_Brat"
$ws.Range("B6").Value = "http://dbpedia.org/resource/NAMUR"
$ws.Range("C6").Value = "'92032"
$ws.Range("D6").Value = "'5100"
$ws.Range("E6").Value = "'56"
$ws.Range("F6").Value = "http://dbpedia.org/resource/DV"
$ws.Range("G6").Value = "'1431"

# --- Row 7 ---
$ws.Range("A7").Value = "http://dbpedia.org/resource/Champion_(sportswear)"
$ws.Range("B7").Value = "http://dbpedia.org/resource/NAMUR"
$ws.Range("C7").Value = "'92024"
$ws.Range("D7").Value = "'5020"
$ws.Range("E7").Value = "'35"
$ws.Range("F7").Value = "http://dbpedia.org/resource/CH"
$ws.Range("G7").Value = "'1717"

# --- Row 8 ---
$ws.Range("A8").Value = "Daussoulx"
$ws.Range("B8").Value = "http://dbpedia.org/resource/NAMUR"
$ws.Range("C8").Value = "'92031"
$ws.Range("D8").Value = "'5020"
$ws.Range("E8").Value = "'18"
$ws.Range("F8").Value = "http://dbpedia.org/resource/DX"
$ws.Range("G8").Value = "'661"

# --- Row 9 ---
$ws.Range("A9").Value = "http://dbpedia.org/resource/Fort_de_Cognelée"
$ws.Range("B9").Value = "http://dbpedia.org/resource/NAMUR"
$ws.Range("C9").Value = "'92025"
$ws.Range("D9").Value = "'5022"
$ws.Range("E9").Value = "'19"
$ws.Range("F9").Value = "http://dbpedia.org/resource/CG"
$ws.Range("G9").Value = "'797"

# --- Row 10 ---
$ws.Range("A10").Value = "Vedrin"
$ws.Range("B10").Value = "http://dbpedia.org/resource/NAMUR"
$ws.Range("C10").Value = "'92128"
$ws.Range("D10").Value = "'5020"
$ws.Range("E10").Value = "'116"
$ws.Range("F10").Value = "http://dbpedia.org/resource/VD"
$ws.Range("G10").Value = "'6580"

# --- Row 11 ---
$ws.Range("A11").Value = "http://dbpedia.org/resource/Sleeze_Beez"
$ws.Range("B11").Value = "http://dbpedia.org/resource/NAMUR"
$ws.Range("C11").Value = "'92009"
$ws.Range("D11").Value = "'5000"
$ws.Range("E11").Value = "'36"
$ws.Range("F11").Value = "http://dbpedia.org/resource/BZ"
$ws.Range("G11").Value = "'1509"

# --- Row 12 ---
$ws.Range("A12").Value = "http://dbpedia.org/resource/Saint-Servais"
$ws.Range("B12").Value = "http://dbpedia.org/resource/NAMUR"
$ws.Range("C12").Value = "'92108"
$ws.Range("D12").Value = "'5002"
$ws.Range("E12").Value = "'87"
$ws.Range("F12").Value = "http://dbpedia.org/resource/SS_(band)"
$ws.Range("G12").Value = "'9299"

# Clear the quote-prefix formatting residue from the forced-text numeric cells,
# so they end up as plain unstyled text cells (matching the target).
$ws.Range("C2:C12").ClearFormats()
$ws.Range("D2:D12").ClearFormats()
$ws.Range("E2:E12").ClearFormats()
$ws.Range("G2:G12").ClearFormats()
